$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = 112241880
$ws.Range("B5").Value = 78713
$ws.Range("E5").Value = 6458
$ws.Range("F5").Value = 'Lunglav'
$ws.Range("G5").Value = 'Lobaria pulmonaria'
$ws.Range("H5").Value = '(L.) Hoffm.'
$ws.Range("Q5").Value = 554164
$ws.Range("R5").Value = 7007925

# Row 6
$ws.Range("B6").Value = 78713

# Row 7
$ws.Range("A7").Value = 112241871
$ws.Range("B7").Value = 89553
$ws.Range("E7").Value = 1202
$ws.Range("F7").Value = 'Ullticka'
$ws.Range("G7").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H7").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q7").Value = 554086
$ws.Range("R7").Value = 7008007

# Row 8
$ws.Range("A8").Value = 112241891
$ws.Range("B8").Value = 78713
$ws.Range("Q8").Value = 553968
$ws.Range("R8").Value = 7008002

# Row 9
$ws.Range("A9").Value = 112241870
$ws.Range("B9").Value = 89553
$ws.Range("E9").Value = 1202
$ws.Range("F9").Value = 'Ullticka'
$ws.Range("G9").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H9").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q9").Value = 554078
$ws.Range("R9").Value = 7008009

# Row 10
$ws.Range("A10").Value = 112241868
$ws.Range("B10").Value = 89553
$ws.Range("E10").Value = 1202
$ws.Range("F10").Value = 'Ullticka'
$ws.Range("G10").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H10").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("I10").Value = ''
$ws.Range("K10").ClearContents()
$ws.Range("L10").ClearContents()
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()
$ws.Range("Q10").Value = 554007
$ws.Range("R10").Value = 7007988

# Row 11
$ws.Range("A11").Value = 112241881
$ws.Range("B11").Value = 78713
$ws.Range("Q11").Value = 554185
$ws.Range("R11").Value = 7007926

# Row 12
$ws.Range("A12").Value = 112241878
$ws.Range("B12").Value = 78713
$ws.Range("E12").Value = 6458
$ws.Range("F12").Value = 'Lunglav'
$ws.Range("G12").Value = 'Lobaria pulmonaria'
$ws.Range("H12").Value = '(L.) Hoffm.'
$ws.Range("Q12").Value = 554158
$ws.Range("R12").Value = 7007927

# Row 13
$ws.Range("A13").Value = 112241869
$ws.Range("B13").Value = 89553
$ws.Range("Q13").Value = 553972
$ws.Range("R13").Value = 7008047

# Row 14
$ws.Range("A14").Value = 112241882
$ws.Range("B14").Value = 78713
$ws.Range("Q14").Value = 553964
$ws.Range("R14").Value = 7008000

# Row 15
$ws.Range("A15").Value = 112241877
$ws.Range("B15").Value = 89993
$ws.Range("D15").Value = 'VU'
$ws.Range("E15").Value = 1209
$ws.Range("F15").Value = 'Rynkskinn'
$ws.Range("G15").Value = 'Phlebia centrifuga'
$ws.Range("H15").Value = 'P.Karst.'
$ws.Range("Q15").Value = 554081
$ws.Range("R15").Value = 7007966

# Row 16
$ws.Range("A16").Value = 112241867
$ws.Range("B16").Value = 89553
$ws.Range("Q16").Value = 554177
$ws.Range("R16").Value = 7007859

# Row 17
$ws.Range("A17").Value = 112241876
$ws.Range("B17").Value = 56446
$ws.Range("D17").Value = 'NT'
$ws.Range("E17").Value = 100049
$ws.Range("F17").Value = 'Spillkråka'
$ws.Range("G17").Value = 'Dryocopus martius'
$ws.Range("H17").Value = '(Linnaeus, 1758)'
$ws.Range("I17").Value = '1'
$ws.Range("K17").Value = ''
$ws.Range("L17").Value = ''
$ws.Range("M17").Value = 'födosökande'
$ws.Range("N17").Value = 'observerad'
$ws.Range("Q17").Value = 554066
$ws.Range("R17").Value = 7008014

# Row 18
$ws.Range("B18").Value = 89499

# Row 19
$ws.Range("B19").Value = 90240
